$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A98").Value = 97
$ws.Range("B98").Value = 1
$ws.Range("C98").Value = "2024-06-16 23:15:27"
$ws.Range("D98").Value = 200
$ws.Range("E98").Value = 6

$ws.Range("A99").Value = 98
$ws.Range("B99").Value = 2
$ws.Range("C99").Value = "2024-06-16 23:15:28"
$ws.Range("D99").Value = 200
$ws.Range("E99").Value = 2
